# Auto-generated Excel COM-interop script
# Applies updated market-board price/profit figures to the Leviathan Profits sheets.
$wb = $excel.ActiveWorkbook

$updates = @{}

$updates["ALC"] = @{
    "H33" = 399.9
    "I33" = 358.16666
    "J33" = 462.5
    "K33" = 358.16666
    "N33" = -920.5
    "H64" = 4668.25
    "I64" = 3437.5
    "J64" = 5899
    "K64" = 3437.5
    "L64" = 5899
    "M64" = -3189.5
    "N64" = -6395
    "H67" = 4668.25
    "I67" = 3437.5
    "J67" = 5899
    "K67" = 3437.5
    "L67" = 5899
    "M67" = -2579.5
    "N67" = -7615
    "H98" = 813.93335
    "I98" = 580.4717000000001
    "K98" = 580.4717000000001
    "M98" = 917.5282999999999
    "H112" = 1440.7894
    "J112" = 1565.7333
    "L112" = 4697.199900000001
    "N112" = -6913.199900000001
    "H122" = 813.93335
    "I122" = 580.4717000000001
    "K122" = 1741.4151
    "M122" = 708.5848999999998
    "H132" = 1068.1052
    "I132" = 1068.1052
    "K132" = 3204.3156
    "M132" = -674.3155999999999
    "H137" = 2012.7812
    "I137" = 1927.4
    "J137" = 2155.0833
    "K137" = 5782.200000000001
    "L137" = 6465.249899999999
    "M137" = -3232.200000000001
    "N137" = -11565.2499
    "H138" = 1796.2927
    "I138" = 1805.3125
    "K138" = 5415.9375
    "M138" = -275.9375
}

$updates["ARM"] = @{
    "H32" = 19220.234
    "I32" = 3953.1714
    "J32" = 90466.53
    "K32" = 3953.1714
    "L32" = 90466.53
    "M32" = -3666.1714
    "N32" = -91040.53
    "H61" = 1716.5
    "I61" = 1763.5385
    "J61" = 1105
    "K61" = 1763.5385
    "L61" = 1105
    "N61" = -1529
    "H63" = 1922.3846
    "I63" = 1922.3846
    "K63" = 1922.3846
    "M63" = -1236.3846
    "H66" = 1922.3846
    "I66" = 1922.3846
    "K66" = 9611.923000000001
    "M66" = -6179.923000000001
    "H74" = 1515.8438
    "I74" = 1258.8966
    "J74" = 3999.6667
    "K74" = 1258.8966
    "L74" = 3999.6667
    "M74" = -384.8966
    "N74" = -5747.6667
    "H77" = 1515.8438
    "I77" = 1258.8966
    "J77" = 3999.6667
    "K77" = 6294.483
    "L77" = 19998.3335
    "M77" = -1926.483
    "N77" = -28734.3335
    "H132" = 1219.8206
    "I132" = 1154.8055
    "K132" = 3464.4165
    "M132" = -934.4164999999998
    "H136" = 1716.5
    "I136" = 1763.5385
    "J136" = 1105
    "K136" = 5290.6155
    "L136" = 3315
    "N136" = -8415
}

$updates["CRP"] = @{
    "H22" = 677.6429000000001
    "I22" = 518.7
    "J22" = 1075
    "K22" = 518.7
    "L22" = 1075
    "M22" = -168.7
    "N22" = -1775
    "H31" = 28105.02
    "I31" = 33056.09
    "J31" = 19935.75
    "K31" = 33056.09
    "L31" = 19935.75
    "M31" = -32761.09
    "N31" = -20525.75
    "H34" = 28105.02
    "I34" = 33056.09
    "J34" = 19935.75
    "K34" = 33056.09
    "L34" = 19935.75
    "M34" = -32854.09
    "N34" = -20339.75
    "H107" = 1792.2727
    "I107" = 1794.3077
    "J107" = 1789.3334
    "K107" = 1794.3077
    "L107" = 1789.3334
    "M107" = 125.6922999999999
    "N107" = -5629.3334
    "H134" = 2151.0986
    "I134" = 2054.9673
    "J134" = 2737.5
    "K134" = 6164.901899999999
    "L134" = 8212.5
    "M134" = -3629.901899999999
    "N134" = -13282.5
}

$updates["CUL"] = @{
    "H37" = 500047500
    "J37" = 500047500
    "L37" = 1500142500
    "N37" = -1500142724
    "H60" = 1230.8
    "I60" = 788.5
    "J60" = 3000
    "K60" = 2365.5
    "L60" = 9000
    "N60" = -9502
    "H131" = 14862.637
    "I131" = 727.8889
    "J131" = 24648.23
    "K131" = 2183.6667
    "L131" = 73944.69
    "M131" = 2856.3333
    "N131" = -84024.69
}

$updates["GSM"] = @{
    "H107" = 55559300
    "I107" = 200
    "J107" = 62504188
    "K107" = 200
    "L107" = 62504188
    "M107" = 1720
    "N107" = -62508028
    "H113" = 2799.1428
    "I113" = 3248.75
    "J113" = 2619.3
    "K113" = 3248.75
    "L113" = 2619.3
    "M113" = -1078.75
    "N113" = -6959.3
}

$updates["LTW"] = @{
    "H98" = 30355
    "J98" = 30355
    "N98" = -36345
    "H132" = 2560.2246
    "I132" = 2192.682
    "J132" = 5794.6
    "K132" = 6578.045999999999
    "L132" = 17383.8
    "M132" = -4048.045999999999
    "N132" = -22443.8
}

$updates["WVR"] = @{
    "H107" = 20005072
    "I107" = 7170.5
    "J107" = 55556896
    "K107" = 21511.5
    "L107" = 166670688
    "M107" = -19591.5
    "N107" = -166674528
    "H132" = 3143.348
    "I132" = 3318.976
    "J132" = 1299.25
    "K132" = 9956.928
    "L132" = 3897.75
    "M132" = -7426.928
    "N132" = -8957.75
    "H136" = 697.1905
    "I136" = 612.6842
    "K136" = 1838.0526
    "M136" = 711.9474
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}

Write-Output "Applied $($updates.Values | ForEach-Object { $_.Count } | Measure-Object -Sum | Select-Object -ExpandProperty Sum) cell updates across $($updates.Keys.Count) sheets."